$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column V (22nd column), shifting existing
# columns V:X to W:Y. The new column inherits formatting from column U.
$ws.Columns("V").Insert()

# Set the header text for the newly inserted column.
$ws.Range("V4").Value = "Фикс сумма"

# Extend the filter database defined name to include the new column.
$wb.Names.Item("_xlnm._FilterDatabase").RefersTo = "=TDSheet!`$A`$4:`$Y`$4"

# Match the author's final cursor position after the edit.
[void]$ws.Range("V5").Select()
